$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48, shifting existing rows 48:128 down to 49:129
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new data record
$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = 44540
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = 100112009
$ws.Range("G48").Value = "Acelga"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 3500
$ws.Range("L48").Value = 3500
$ws.Range("M48").Value = 3500
$ws.Range("N48").Value = "$/docena de atados (4 kilos)"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 875
$ws.Range("Q48").Value = 4
$ws.Range("R48").Value = "Hortaliza"

# Match the D-column number format used by the rest of the date column (style index 2)
$ws.Range("D48").NumberFormat = $ws.Range("D49").NumberFormat
